# Commit 02: Team 10 Finish 2nd Progress
# - Mark the "Long" / "Hien" team-member assignments and "Xong" (done)
#   status for the tasks in groups 2 and 3, centering the status column.
# - Update the saved scroll position / selection of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectPlan")

# Rows 12-16 -> Long finished his part; rows 17-20 -> Hien finished hers.
$longRows1 = 12, 13, 14, 15, 16
$hienRows1 = 17, 18, 19, 20
foreach ($r in $longRows1) {
    $ws.Cells.Item($r, 5).Value = "Long"
}
foreach ($r in $hienRows1) {
    $ws.Cells.Item($r, 5).Value = "Hien"
}
foreach ($r in (12..20)) {
    $ws.Cells.Item($r, 6).Value = "Xong"
}
$ws.Range("F12:F20").HorizontalAlignment = -4108

# Rows 22-25 -> Long finished his part; rows 26-29 -> Hien finished hers.
$longRows2 = 22, 23, 24, 25
$hienRows2 = 26, 27, 28, 29
foreach ($r in $longRows2) {
    $ws.Cells.Item($r, 5).Value = "Long"
}
foreach ($r in $hienRows2) {
    $ws.Cells.Item($r, 5).Value = "Hien"
}
foreach ($r in (22..29)) {
    $ws.Cells.Item($r, 6).Value = "Xong"
}
$ws.Range("F22:F29").HorizontalAlignment = -4108

# Scroll the sheet down and move the active selection, matching the
# author's view when they saved the workbook.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("H30").Select()
